$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J9").Value = 3502
$ws.Range("I9").Value = 365.33334
$ws.Range("M9").Value = -196.33334
$ws.Range("K9").Value = 365.33334
$ws.Range("L9").Value = 3502
$ws.Range("N9").Value = -3840
$ws.Range("H9").Value = 679
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("L23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 1934.4
$ws.Range("J29").Value = 2646.4443
$ws.Range("L29").Value = 7939.3329
$ws.Range("N29").Value = -8501.332900000001
$ws.Range("J64").Value = 7850
$ws.Range("L64").Value = 7850
$ws.Range("H64").Value = 6825
$ws.Range("N64").Value = -8346
$ws.Range("L67").Value = 7850
$ws.Range("N67").Value = -9566
$ws.Range("H67").Value = 6825
$ws.Range("J67").Value = 7850
$ws.Range("N70").Value = -12087.75
$ws.Range("M70").Value = -7830
$ws.Range("H70").Value = 3117.9092
$ws.Range("J70").Value = 3849.25
$ws.Range("I70").Value = 2700
$ws.Range("K70").Value = 8100
$ws.Range("L70").Value = 11547.75
$ws.Range("I73").Value = 2700
$ws.Range("J73").Value = 3849.25
$ws.Range("K73").Value = 8100
$ws.Range("M73").Value = -7164
$ws.Range("L73").Value = 11547.75
$ws.Range("H73").Value = 3117.9092
$ws.Range("N73").Value = -13419.75
$ws.Range("H96").Value = 1705.1666
$ws.Range("I96").Value = 1757.3077
$ws.Range("K96").Value = 5271.9231
$ws.Range("J96").Value = 1569.6
$ws.Range("M96").Value = -3898.9231
$ws.Range("L96").Value = 4708.799999999999
$ws.Range("N96").Value = -7454.799999999999
$ws.Range("J100").Value = 4750
$ws.Range("I100").Value = 2468.7144
$ws.Range("K100").Value = 2468.7144
$ws.Range("M100").Value = -1927.7144
$ws.Range("L100").Value = 4750
$ws.Range("N100").Value = -5832
$ws.Range("H100").Value = 2975.6667
$ws.Range("M106").Value = -3318.5833
$ws.Range("L106").Value = 3999.3333
$ws.Range("H106").Value = 3959.5334
$ws.Range("N106").Value = -5261.3333
$ws.Range("J106").Value = 3999.3333
$ws.Range("I106").Value = 3949.5833
$ws.Range("K106").Value = 3949.5833
$ws.Range("M137").Value = -3538.9092
$ws.Range("H137").Value = 8876.079
$ws.Range("J137").Value = 16407.166
$ws.Range("I137").Value = 2029.6364
$ws.Range("K137").Value = 6088.9092
$ws.Range("L137").Value = 49221.49800000001
$ws.Range("N137").Value = -54321.49800000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 50000
$ws.Range("J23").Value = 50000
$ws.Range("L23").Value = 50000
$ws.Range("N23").Value = -50518
$ws.Range("J37").Value = 30000
$ws.Range("I37").Value = 10000
$ws.Range("K37").Value = 10000
$ws.Range("M37").Value = -9727
$ws.Range("L37").Value = 30000
$ws.Range("H37").Value = 20000
$ws.Range("N37").Value = -30546
$ws.Range("I45").Value = 1368.8889
$ws.Range("K45").Value = 1368.8889
$ws.Range("J45").Value = 1926.125
$ws.Range("M45").Value = -991.8888999999999
$ws.Range("L45").Value = 1926.125
$ws.Range("H45").Value = 1631.1177
$ws.Range("N45").Value = -2680.125
$ws.Range("I74").Value = 3996.3125
$ws.Range("M74").Value = -3122.3125
$ws.Range("K74").Value = 3996.3125
$ws.Range("H74").Value = 15112.81
$ws.Range("I77").Value = 3996.3125
$ws.Range("K77").Value = 19981.5625
$ws.Range("M77").Value = -15613.5625
$ws.Range("H77").Value = 15112.81
$ws.Range("I97").Value = 3395.2727
$ws.Range("K97").Value = 3395.2727
$ws.Range("M97").Value = -2899.2727
$ws.Range("H97").Value = 3074.5386
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 29235.28
$ws.Range("N20").Value = -40560.855
$ws.Range("I20").Value = 18896.045
$ws.Range("J20").Value = 40066.855
$ws.Range("K20").Value = 18896.045
$ws.Range("M20").Value = -18649.045
$ws.Range("L20").Value = 40066.855
$ws.Range("M86").Value = -3258.8
$ws.Range("H86").Value = 4484.8335
$ws.Range("I86").Value = 4381.8
$ws.Range("K86").Value = 4381.8
$ws.Range("I89").Value = 4381.8
$ws.Range("K89").Value = 21909
$ws.Range("M89").Value = -16293
$ws.Range("H89").Value = 4484.8335
$ws.Range("I134").Value = 14093.8
$ws.Range("J134").Value = 43605.6
$ws.Range("K134").Value = 42281.39999999999
$ws.Range("N134").Value = -135886.8
$ws.Range("M134").Value = -39746.39999999999
$ws.Range("L134").Value = 130816.8
$ws.Range("H134").Value = 19996.16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 47276.945
$ws.Range("N31").Value = -24155
$ws.Range("I31").Value = 60527.734
$ws.Range("J31").Value = 23565
$ws.Range("K31").Value = 60527.734
$ws.Range("M31").Value = -60232.734
$ws.Range("L31").Value = 23565
$ws.Range("I34").Value = 60527.734
$ws.Range("J34").Value = 23565
$ws.Range("K34").Value = 60527.734
$ws.Range("M34").Value = -60325.734
$ws.Range("L34").Value = 23565
$ws.Range("N34").Value = -23969
$ws.Range("H34").Value = 47276.945
$ws.Range("I58").Value = 10512.429
$ws.Range("J58").Value = 32546.363
$ws.Range("K58").Value = 10512.429
$ws.Range("M58").Value = -10309.429
$ws.Range("L58").Value = 32546.363
$ws.Range("H58").Value = 23977.611
$ws.Range("N58").Value = -32952.363
$ws.Range("J59").Value = 5000
$ws.Range("N59").Value = -7290
$ws.Range("L59").Value = 5000
$ws.Range("H59").Value = 5000
$ws.Range("J122").Value = 2856.5715
$ws.Range("I122").Value = 2637.1667
$ws.Range("K122").Value = 7911.500100000001
$ws.Range("L122").Value = 8569.7145
$ws.Range("N122").Value = -13469.7145
$ws.Range("M122").Value = -5461.500100000001
$ws.Range("H122").Value = 2755.3076
$ws.Range("J136").Value = 32546.363
$ws.Range("I136").Value = 10512.429
$ws.Range("M136").Value = -28987.287
$ws.Range("K136").Value = 31537.287
$ws.Range("L136").Value = 97639.08900000001
$ws.Range("N136").Value = -102739.089
$ws.Range("H136").Value = 23977.611
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 258.75
$ws.Range("K5").Value = 776.25
$ws.Range("M5").Value = -664.25
$ws.Range("H5").Value = 435.13333
$ws.Range("H129").Value = 2105.276
$ws.Range("J129").Value = 3578.3333
$ws.Range("I129").Value = 527
$ws.Range("K129").Value = 1581
$ws.Range("L129").Value = 10734.9999
$ws.Range("N129").Value = -20734.9999
$ws.Range("M129").Value = 3419
$ws.Range("I131").Value = 923.625
$ws.Range("J131").Value = 1470.1195
$ws.Range("K131").Value = 2770.875
$ws.Range("H131").Value = 1426.4
$ws.Range("M131").Value = 2269.125
$ws.Range("L131").Value = 4410.3585
$ws.Range("N131").Value = -14490.3585
$ws.Range("I132").Value = 1524.1666
$ws.Range("K132").Value = 13717.4994
$ws.Range("H132").Value = 1516715.8
$ws.Range("M132").Value = -11187.4994
$ws.Range("M135").Value = 206.25
$ws.Range("H135").Value = 435.13333
$ws.Range("K135").Value = 2328.75
$ws.Range("I135").Value = 258.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L107").Value = 726.5
$ws.Range("H107").Value = 613.9231
$ws.Range("N107").Value = -4566.5
$ws.Range("J107").Value = 726.5
$ws.Range("J121").Value = 99999
$ws.Range("L121").Value = 99999
$ws.Range("N121").Value = -103493
$ws.Range("H121").Value = 99999
$ws.Range("I132").Value = 5643.2354
$ws.Range("J132").Value = 745098.6
$ws.Range("K132").Value = 16929.7062
$ws.Range("L132").Value = 2235295.8
$ws.Range("H132").Value = 430830.1
$ws.Range("N132").Value = -2240355.8
$ws.Range("M132").Value = -14399.7062
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M93").Value = -8981.5
$ws.Range("L93").Value = 3997.5
$ws.Range("H93").Value = 8671.5
$ws.Range("N93").Value = -6493.5
$ws.Range("I93").Value = 10229.5
$ws.Range("J93").Value = 3997.5
$ws.Range("K93").Value = 10229.5
$ws.Range("J122").Value = 7112.625
$ws.Range("I122").Value = 11415.857
$ws.Range("K122").Value = 34247.571
$ws.Range("L122").Value = 21337.875
$ws.Range("N122").Value = -26237.875
$ws.Range("M122").Value = -31797.571
$ws.Range("H122").Value = 8422.305
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I122").Value = 2925.56
$ws.Range("K122").Value = 8776.68
$ws.Range("M122").Value = -6326.68
$ws.Range("H122").Value = 4192.95
$ws.Range("I132").Value = 3304.2
$ws.Range("J132").Value = 5330283
$ws.Range("K132").Value = 9912.599999999999
$ws.Range("L132").Value = 15990849
$ws.Range("H132").Value = 2793626.5
$ws.Range("N132").Value = -15995909
$ws.Range("M132").Value = -7382.599999999999

Write-Host "Done applying edits"